$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sheet" (first worksheet): update release date, replace the
# FH-SM05R/optics accessory row with the new FHV7H-C016-C / K50RPLPGREQP
# entry, and drop the now-unused trailing rows.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Release date (E3) - keep it as literal text (not an auto-converted date
# serial) and keep the original cell style by pasting the format back from
# the sibling cell E2.
$ws1.Range("E3").NumberFormat = "@"
$ws1.Range("E3").Value = "07.08.2025"
$ws1.Range("E2").Copy()
$ws1.Range("E3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 6 content: new camera / optics accessory part numbers.
$ws1.Range("C6").Value = "FHV7H-C016-C"
$ws1.Range("D6").Value = ""
$ws1.Range("E6").Value = "K50RPLPGREQP"
$ws1.Range("F6").Value = "Alternativa: FZ-LEH50"
$ws1.Rows.Item(6).RowHeight = 15

# Row 7 becomes the blank footer-styled row (formerly row 9): copy its
# format+content over, which also releases the merges touching old rows
# 6-8, then clear any leftover text.
$ws1.Range("A9:J9").Copy($ws1.Range("A7:J7"))
$excel.CutCopyMode = $false
$ws1.Range("A7:J7").ClearContents()

# Rows 8 and 9 are no longer needed.
$ws1.Rows.Item(9).Delete()
$ws1.Rows.Item(8).Delete()

# ---------------------------------------------------------------------------
# Sheet "Kusovník" (BOM sheet): same part-number swap, quantities adjusted,
# and the now-removed controller-accessory entry (row 4) dropped.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A3").Value = "FHV7H-C016-C"
$ws2.Range("B3").Value = 1
$ws2.Range("D3").Value = "K50RPLPGREQP"

# Drop the "Příslušenství" entries entirely (M3:N3 content, plus the whole
# M4:N4 row) - Clear() (not just ClearContents) removes the cells outright
# instead of leaving blank placeholders, so row 4 disappears and the sheet
# dimension shrinks back to N3.
$ws2.Range("M3:N4").Clear()

Write-Host "edit complete"
